$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.372.71'
$ws.Range("E2").Value = '  -0.36%  '

$ws.Range("D3").Value = '3.071.05'
$ws.Range("E3").Value = '  +1.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '395.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.534'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.82%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.585'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.56%  '

$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0850'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.08%  '

$ws.Range("D13").Value = '3.554.29'
$ws.Range("E13").Value = '  +1.22%  '

$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("E15").Value = '  -1.02%  '

$ws.Range("D16").Value = '3.076.23'
$ws.Range("E16").Value = '  +1.39%  '

$ws.Range("E17").Value = '  +4.84%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.53%  '

$ws.Range("D19").Value = '51.364.90'
$ws.Range("E19").Value = '  -0.42%  '

$ws.Range("E20").Value = '  +2.46%  '

$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("D22").Value = '0.0₃0957'
$ws.Range("E22").Value = '  -0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("E24").Value = '  -0.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.23%  '

$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.21%  '

$ws.Range("E30").Value = '  -6.73%  '

$ws.Range("E31").Value = '  -1.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0486'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '35.95'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.52%  '

$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '50.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("E38").Value = '  -1.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.291'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.85%  '

$ws.Range("E40").Value = '  +6.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '128.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.67'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.73%  '

$ws.Range("E43").Value = '  -1.06%  '

$ws.Range("E44").Value = '  -1.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.52'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.35%  '

$ws.Range("E48").Value = '  -2.55%  '

$ws.Range("D49").Value = '2.073.34'
$ws.Range("E49").Value = '  +1.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.891'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.95%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.515'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.12%  '
